# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the da7b42e3-... item, in both the zh-cn and de-de sheets,
# plus the rolled-up status on the Overview sheet.

$wb = $excel.ActiveWorkbook

$zhError = "Handback file name: 3dx5yv20.q0m is different with handoff file name: da7b42e3-9a72-44b6-b4dd-8d4db54fe8d5.0f58b8f6886e80f26f281c1099dcac32c9c8b94e.zh-cn."
$deError = "Handback file name: 3dx5yv20.q0m is different with handoff file name: da7b42e3-9a72-44b6-b4dd-8d4db54fe8d5.0f58b8f6886e80f26f281c1099dcac32c9c8b94e.de-de."
$newStatus = "Handback transform failed"

# --- Overview sheet: roll up status changes for the da7b42e3-... row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# NOTE: the .xlsx <col width=.../> attribute stores character-width + 5/6
# (the default-font padding), so to land on a stored width of exactly 40
# we must request ColumnWidth = 40 - 5/6.
$targetColWidth = 235 / 6

# --- zh-cn sheet: set Status + Error Detail for row 3, widen column P ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = $targetColWidth

# --- de-de sheet: set Status + Error Detail for row 3, widen column P ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = $targetColWidth
